$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5357409715652466
$ws.Range("B1").Value = 0.6755579113960266
$ws.Range("C1").Value = 0.9819005727767944
$ws.Range("D1").Value = 3.811616659164429
$ws.Range("E1").Value = 5.695967197418213
